$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the swapped rows 10 and 11 (Polygon / Chainlink reordered)
$ws.Range("A10").Value = "Polygon"
$ws.Range("B10").Value = "MATIC-USD"
$ws.Range("A11").Value = "Chainlink"
$ws.Range("B11").Value = "LINK-USD"

# Update Market Cap (column C) values for rows 2-26
$ws.Range("C2").Value = 737133302290.2677
$ws.Range("C3").Value = 247244493778.9961
$ws.Range("C4").Value = 38390659845.19044
$ws.Range("C5").Value = 34685910398.30289
$ws.Range("C6").Value = 27532074178.23908
$ws.Range("C7").Value = 13393577812.97504
$ws.Range("C8").Value = 10819894707.74953
$ws.Range("C9").Value = 9242314400.824749
$ws.Range("C10").Value = 8507732893.082836
$ws.Range("C11").Value = 8325744672.079366
$ws.Range("C12").Value = 8189957692.931383
$ws.Range("C13").Value = 7354641270.491064
$ws.Range("C14").Value = 6982621694.335486
$ws.Range("C15").Value = 6178344669.330246
$ws.Range("C16").Value = 5463567574.120149
$ws.Range("C17").Value = 5176763885.483999
$ws.Range("C18").Value = 4682341026.782693
$ws.Range("C19").Value = 3777671236.637238
$ws.Range("C20").Value = 3648354300.80387
$ws.Range("C21").Value = 3614419588.569764
$ws.Range("C22").Value = 3398722194.285174
$ws.Range("C23").Value = 2979199919.429708
$ws.Range("C24").Value = 2856082645.519599
$ws.Range("C25").Value = 2624592537.812566
$ws.Range("C26").Value = 2522164756.451115
